# Refresh the crypto price ("Price", column D) and hourly volume
# change ("Volume(1h)", column E) figures in rows 2-51 of Sheet1,
# as produced by the scheduled GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.877.37'
$ws.Range('E2').Value = '  -0.77%  '
$ws.Range('D3').Value = '1.620.19'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('E4').Value = '  -0.57%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '212.99'
$cell.ClearFormats()
$ws.Range('E5').Value = '  -1.39%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '0.500'
$cell.ClearFormats()
$ws.Range('E6').Value = '  -1.46%  '
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('E8').Value = '  -1.45%  '
$ws.Range('E9').Value = '  -3.15%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '18.49'
$cell.ClearFormats()
$ws.Range('E10').Value = '  -5.01%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.0790'
$cell.ClearFormats()
$ws.Range('E11').Value = '  -0.83%  '
$ws.Range('D12').Value = '1.845.44'
$ws.Range('E12').Value = '  -1.45%  '
$ws.Range('D13').Value = '1.614.38'
$ws.Range('E13').Value = '  -2.93%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '4.15'
$cell.ClearFormats()
$ws.Range('E14').Value = '  -2.52%  '
$ws.Range('E15').Value = '  -3.19%  '
$ws.Range('D16').Value = '25.887.79'
$ws.Range('E16').Value = '  -0.75%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '61.57'
$cell.ClearFormats()
$ws.Range('E17').Value = '  -2.87%  '
$ws.Range('E18').Value = '  -2.87%  '
$ws.Range('E19').Value = '  -0.39%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '191.84'
$cell.ClearFormats()
$ws.Range('E20').Value = '  -1.20%  '
$ws.Range('E21').Value = '  -2.04%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '9.50'
$cell.ClearFormats()
$ws.Range('E22').Value = '  -2.55%  '
$ws.Range('E23').Value = '  -2.33%  '
$ws.Range('E24').Value = '  +2.75%  '
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('E26').Value = '  -0.58%  '
$ws.Range('E27').Value = '  -4.02%  '
$ws.Range('E28').Value = '  -2.70%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '15.19'
$cell.ClearFormats()
$ws.Range('E29').Value = '  -2.21%  '
$ws.Range('E30').Value = '  -1.27%  '
$ws.Range('E31').Value = '  -2.18%  '
$ws.Range('E32').Value = '  -3.73%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '3.10'
$cell.ClearFormats()
$ws.Range('E33').Value = '  -5.25%  '
$ws.Range('E34').Value = '  -2.37%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '1.49'
$cell.ClearFormats()
$ws.Range('E35').Value = '  -2.74%  '
$ws.Range('D36').Value = '1.125.92'
$ws.Range('E36').Value = '  -0.37%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '0.840'
$cell.ClearFormats()
$ws.Range('E37').Value = '  -6.71%  '
$ws.Range('E38').Value = '  -4.17%  '
$ws.Range('E39').Value = '  -2.09%  '
$ws.Range('E40').Value = '  -4.00%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '98.20'
$cell.ClearFormats()
$ws.Range('E41').Value = '  -0.38%  '
$ws.Range('D42').Value = '1.755.77'
$ws.Range('E42').Value = '  -1.31%  '
$ws.Range('E43').Value = '  -5.88%  '
$ws.Range('E44').Value = '  -5.45%  '
$ws.Range('E45').Value = '  -1.21%  '
$ws.Range('E46').Value = '  +1.51%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '54.06'
$cell.ClearFormats()
$ws.Range('E47').Value = '  -3.85%  '
$ws.Range('E48').Value = '  -0.44%  '
$ws.Range('E49').Value = '  -1.24%  '
$ws.Range('E50').Value = '  -0.42%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '7.46'
$cell.ClearFormats()
$ws.Range('E51').Value = '  -3.75%  '
